$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "25.820.79"
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  -5.18%  "
$c.ClearFormats()

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.813.12"
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  -4.30%  "
$c.ClearFormats()

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9991"
$c.ClearFormats()

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "276.33"
$c.ClearFormats()

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  -0.30%  "
$c.ClearFormats()

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5090"
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  -5.57%  "
$c.ClearFormats()

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3522"
$c.ClearFormats()
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  -7.17%  "
$c.ClearFormats()

$c = $ws.Range("B9")
$c.NumberFormat = "@"
$c.Value = "Dogecoin"
$c.ClearFormats()
$c = $ws.Range("C9")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$c.ClearFormats()
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06662"
$c.ClearFormats()
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  -8.46%  "
$c.ClearFormats()

$c = $ws.Range("B10")
$c.NumberFormat = "@"
$c.Value = "Solana"
$c.ClearFormats()
$c = $ws.Range("C10")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$c.ClearFormats()
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "20.09"
$c.ClearFormats()
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  -8.69%  "
$c.ClearFormats()

$c = $ws.Range("B11")
$c.NumberFormat = "@"
$c.Value = "Polygon"
$c.ClearFormats()
$c = $ws.Range("C11")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c.ClearFormats()
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.8359"
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  -6.92%  "
$c.ClearFormats()

$c = $ws.Range("B12")
$c.NumberFormat = "@"
$c.Value = "TRON"
$c.ClearFormats()
$c = $ws.Range("C12")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c.ClearFormats()
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07816"
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  -4.58%  "
$c.ClearFormats()

$c = $ws.Range("B13")
$c.NumberFormat = "@"
$c.Value = "WrappedEther"
$c.ClearFormats()
$c = $ws.Range("C13")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c.ClearFormats()
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.801.94"
$c.ClearFormats()
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +0.56%  "
$c.ClearFormats()

$c = $ws.Range("B14")
$c.NumberFormat = "@"
$c.Value = "Polkadot"
$c.ClearFormats()
$c = $ws.Range("C14")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c.ClearFormats()
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.076"
$c.ClearFormats()
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -4.98%  "
$c.ClearFormats()

$c = $ws.Range("B15")
$c.NumberFormat = "@"
$c.Value = "Litecoin"
$c.ClearFormats()
$c = $ws.Range("C15")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c.ClearFormats()
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "87.94"
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  -7.21%  "
$c.ClearFormats()

$c = $ws.Range("B16")
$c.NumberFormat = "@"
$c.Value = "BinanceUSD"
$c.ClearFormats()
$c = $ws.Range("C16")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c.ClearFormats()
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.9983"
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  -0.35%  "
$c.ClearFormats()

$c = $ws.Range("B17")
$c.NumberFormat = "@"
$c.Value = "Avalanche"
$c.ClearFormats()
$c = $ws.Range("C17")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c.ClearFormats()
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "13.89"
$c.ClearFormats()
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  -6.45%  "
$c.ClearFormats()

$c = $ws.Range("B18")
$c.NumberFormat = "@"
$c.Value = "ShibaInu"
$c.ClearFormats()
$c = $ws.Range("C18")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c.ClearFormats()
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008024"
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  -7.06%  "
$c.ClearFormats()

$c = $ws.Range("B19")
$c.NumberFormat = "@"
$c.Value = "Dai"
$c.ClearFormats()
$c = $ws.Range("C19")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c.ClearFormats()
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.9988"
$c.ClearFormats()
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  -0.36%  "
$c.ClearFormats()

$c = $ws.Range("B20")
$c.NumberFormat = "@"
$c.Value = "WrappedBTC"
$c.ClearFormats()
$c = $ws.Range("C20")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$c.ClearFormats()
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "25.878.55"
$c.ClearFormats()
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  -4.25%  "
$c.ClearFormats()

$c = $ws.Range("B21")
$c.NumberFormat = "@"
$c.Value = "Uniswap"
$c.ClearFormats()
$c = $ws.Range("C21")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c.ClearFormats()
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.734"
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  -5.95%  "
$c.ClearFormats()

$c = $ws.Range("B22")
$c.NumberFormat = "@"
$c.Value = "Cosmos"
$c.ClearFormats()
$c = $ws.Range("C22")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c.ClearFormats()
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "10.01"
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  -7.38%  "
$c.ClearFormats()

$c = $ws.Range("B23")
$c.NumberFormat = "@"
$c.Value = "Chainlink"
$c.ClearFormats()
$c = $ws.Range("C23")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c.ClearFormats()
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "6.061"
$c.ClearFormats()
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  -6.20%  "
$c.ClearFormats()

$c = $ws.Range("B24")
$c.NumberFormat = "@"
$c.Value = "LidoDAOToken"
$c.ClearFormats()
$c = $ws.Range("C24")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c.ClearFormats()
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.214"
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  -3.44%  "
$c.ClearFormats()

$c = $ws.Range("B25")
$c.NumberFormat = "@"
$c.Value = "Monero"
$c.ClearFormats()
$c = $ws.Range("C25")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c.ClearFormats()
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "141.77"
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  -4.69%  "
$c.ClearFormats()

$c = $ws.Range("B26")
$c.NumberFormat = "@"
$c.Value = "Toncoin"
$c.ClearFormats()
$c = $ws.Range("C26")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c.ClearFormats()
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.655"
$c.ClearFormats()
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  -5.56%  "
$c.ClearFormats()

$c = $ws.Range("B27")
$c.NumberFormat = "@"
$c.Value = "EthereumClassic"
$c.ClearFormats()
$c = $ws.Range("C27")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c.ClearFormats()
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "17.04"
$c.ClearFormats()
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -7.13%  "
$c.ClearFormats()

$c = $ws.Range("B28")
$c.NumberFormat = "@"
$c.Value = "BitcoinCash"
$c.ClearFormats()
$c = $ws.Range("C28")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c.ClearFormats()
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "108.89"
$c.ClearFormats()
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  -6.30%  "
$c.ClearFormats()

$c = $ws.Range("B29")
$c.NumberFormat = "@"
$c.Value = "InternetComputer(DFINITY)"
$c.ClearFormats()
$c = $ws.Range("C29")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c.ClearFormats()
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "4.339"
$c.ClearFormats()
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  -9.93%  "
$c.ClearFormats()

$c = $ws.Range("B30")
$c.NumberFormat = "@"
$c.Value = "Filecoin"
$c.ClearFormats()
$c = $ws.Range("C30")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c.ClearFormats()
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.219"
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -9.01%  "
$c.ClearFormats()

$c = $ws.Range("B31")
$c.NumberFormat = "@"
$c.Value = "Stellar"
$c.ClearFormats()
$c = $ws.Range("C31")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c.ClearFormats()
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08781"
$c.ClearFormats()
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -4.15%  "
$c.ClearFormats()

$c = $ws.Range("B32")
$c.NumberFormat = "@"
$c.Value = "Hedera"
$c.ClearFormats()
$c = $ws.Range("C32")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c.ClearFormats()
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.04881"
$c.ClearFormats()
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  -3.34%  "
$c.ClearFormats()

$c = $ws.Range("B33")
$c.NumberFormat = "@"
$c.Value = "ImmutableX"
$c.ClearFormats()
$c = $ws.Range("C33")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c.ClearFormats()
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.7338"
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  -10.51%  "
$c.ClearFormats()

$c = $ws.Range("B34")
$c.NumberFormat = "@"
$c.Value = "ARBITRUM"
$c.ClearFormats()
$c = $ws.Range("C34")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c.ClearFormats()
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.140"
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  -6.30%  "
$c.ClearFormats()

$c = $ws.Range("B35")
$c.NumberFormat = "@"
$c.Value = "HuobiToken"
$c.ClearFormats()
$c = $ws.Range("C35")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c.ClearFormats()
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.895"
$c.ClearFormats()
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  -4.16%  "
$c.ClearFormats()

$c = $ws.Range("B36")
$c.NumberFormat = "@"
$c.Value = "Frax"
$c.ClearFormats()
$c = $ws.Range("C36")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$c.ClearFormats()
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.9983"
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  -0.66%  "
$c.ClearFormats()

$c = $ws.Range("B37")
$c.NumberFormat = "@"
$c.Value = "MXToken"
$c.ClearFormats()
$c = $ws.Range("C37")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c.ClearFormats()
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.038"
$c.ClearFormats()
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  -7.56%  "
$c.ClearFormats()

$c = $ws.Range("B38")
$c.NumberFormat = "@"
$c.Value = "TheSandbox"
$c.ClearFormats()
$c = $ws.Range("C38")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c.ClearFormats()
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.5240"
$c.ClearFormats()
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -12.24%  "
$c.ClearFormats()

$c = $ws.Range("B39")
$c.NumberFormat = "@"
$c.Value = "VeChain"
$c.ClearFormats()
$c = $ws.Range("C39")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c.ClearFormats()
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01856"
$c.ClearFormats()
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -6.58%  "
$c.ClearFormats()

$c = $ws.Range("B40")
$c.NumberFormat = "@"
$c.Value = "RenderToken"
$c.ClearFormats()
$c = $ws.Range("C40")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c.ClearFormats()
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.296"
$c.ClearFormats()
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  -14.12%  "
$c.ClearFormats()

$c = $ws.Range("B41")
$c.NumberFormat = "@"
$c.Value = "TrustWalletToken"
$c.ClearFormats()
$c = $ws.Range("C41")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c.ClearFormats()
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.9529"
$c.ClearFormats()
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  -11.34%  "
$c.ClearFormats()

$c = $ws.Range("B42")
$c.NumberFormat = "@"
$c.Value = "Quant"
$c.ClearFormats()
$c = $ws.Range("C42")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c.ClearFormats()
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "112.81"
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  -1.84%  "
$c.ClearFormats()

$c = $ws.Range("B43")
$c.NumberFormat = "@"
$c.Value = "FraxShare"
$c.ClearFormats()
$c = $ws.Range("C43")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c.ClearFormats()
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "6.193"
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -6.55%  "
$c.ClearFormats()

$c = $ws.Range("B44")
$c.NumberFormat = "@"
$c.Value = "Aptos"
$c.ClearFormats()
$c = $ws.Range("C44")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c.ClearFormats()
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.119"
$c.ClearFormats()
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  -11.86%  "
$c.ClearFormats()

$c = $ws.Range("B45")
$c.NumberFormat = "@"
$c.Value = "PaxDollar"
$c.ClearFormats()
$c = $ws.Range("C45")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c.ClearFormats()
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.9984"
$c.ClearFormats()
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  -0.32%  "
$c.ClearFormats()

$c = $ws.Range("B46")
$c.NumberFormat = "@"
$c.Value = "Decentraland"
$c.ClearFormats()
$c = $ws.Range("C46")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c.ClearFormats()
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.4573"
$c.ClearFormats()
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  -10.32%  "
$c.ClearFormats()

$c = $ws.Range("B47")
$c.NumberFormat = "@"
$c.Value = "Algorand"
$c.ClearFormats()
$c = $ws.Range("C47")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c.ClearFormats()
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.1379"
$c.ClearFormats()
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -9.62%  "
$c.ClearFormats()

$c = $ws.Range("B48")
$c.NumberFormat = "@"
$c.Value = "EnergySwap"
$c.ClearFormats()
$c = $ws.Range("C48")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c.ClearFormats()
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "9.346"
$c.ClearFormats()
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  -8.30%  "
$c.ClearFormats()

$c = $ws.Range("B49")
$c.NumberFormat = "@"
$c.Value = "Elrond"
$c.ClearFormats()
$c = $ws.Range("C49")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$c.ClearFormats()
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "36.22"
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  -4.60%  "
$c.ClearFormats()

$c = $ws.Range("B50")
$c.NumberFormat = "@"
$c.Value = "NEARProtocol"
$c.ClearFormats()
$c = $ws.Range("C50")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c.ClearFormats()
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.500"
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  -7.86%  "
$c.ClearFormats()

$c = $ws.Range("B51")
$c.NumberFormat = "@"
$c.Value = "Cronos"
$c.ClearFormats()
$c = $ws.Range("C51")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c.ClearFormats()
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.05801"
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  -4.70%  "
$c.ClearFormats()
